# Trade #83 closed at 2026-02-17 08:59:23 - unknown UNKNOWN +0.000%
#
# Helper: write a string value into a cell while preventing Excel's
# automatic "looks like a date" (or other) type coercion, and without
# leaving a lingering explicit number-format style on the cell.
function Set-TextValue($cell, [string]$val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.52   # Current Capital
$summary.Range("B4").Value = 0.53      # Total P&L $
$summary.Range("B6").Value = 83        # Total Trades
$summary.Range("B7").Value = 36        # Winning Trades
$summary.Range("B9").Value = 43.37     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet (MarketMaking row)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.52     # Capital
$status.Range("D4").Value = 83         # Trades
$status.Range("E4").Value = 0.53       # P&L $
$status.Range("F4").Value = 0.52       # P&L %
$status.Range("G4").Value = 43.37      # Win Rate %

# ---------------------------------------------------------------
# New trade row (#83) appended to both "All Trades" and
# "MarketMaking" sheets as row 84.
# ---------------------------------------------------------------
$newRow = @{
    A = 83
    B = "2026-02-17"
    C = "08:59:16"
    D = "MarketMaking"
    E = "DOWN"
    F = 0.96
    G = 0.97
    H = "CLOSED"
    I = 1.0417
    J = 0.01
    K = 100.52
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A84").Value = $newRow.A
    Set-TextValue $ws.Range("B84") $newRow.B
    Set-TextValue $ws.Range("C84") $newRow.C
    Set-TextValue $ws.Range("D84") $newRow.D
    Set-TextValue $ws.Range("E84") $newRow.E
    $ws.Range("F84").Value = $newRow.F
    $ws.Range("G84").Value = $newRow.G
    Set-TextValue $ws.Range("H84") $newRow.H
    $ws.Range("I84").Value = $newRow.I
    $ws.Range("J84").Value = $newRow.J
    $ws.Range("K84").Value = $newRow.K
    $ws.Range("L84").Value = $newRow.L
    $ws.Range("M84").Value = $newRow.M
    $ws.Range("N84").Value = $newRow.N
    Set-TextValue $ws.Range("O84") $newRow.O
    Set-TextValue $ws.Range("P84") $newRow.P
    $ws.Range("Q84").Value = $newRow.Q
}
